# Add a "Result" column (C) to the Countries sheet, marking every row as "PASS".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Countries")

$ws.Range("C1").Value = "Result"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = "PASS"
}
